$d = $word.ActiveDocument

# Change 1: "Minimization of data loss" -> "Dealing with data loss"
$d.Content.Find.Execute("Minimization of data loss", $true, $false, $false, $false, $false, $true, 1, $false, "Dealing with data loss", 2)

